# Styling and data updates
# - Adds an emoji (column F) for each "Gluten Free Cereals & Grain" row
#   (rows 134-168 on the "Data" sheet). Column G's formula recomputes
#   automatically to embed the new emoji in its generated string.
# - F159 picks up the plain border style (matching its neighbours) instead
#   of the "Wheat Pasta" row's special highlighted style.
# - Leaves the active selection on F168, matching where the edit ended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# F159 currently carries the row's highlighted style (s="5"); the edited
# workbook shows it switched to the regular bordered style used by the
# rest of the column (s="2"), like its neighbour F158.
$ws.Range("F158").Copy()
$ws.Range("F159").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$emoji = @{
    134 = "🍞"; 135 = "🌾"; 136 = "🍫"; 137 = "🍫"; 138 = "🌽"; 139 = "🌽";
    140 = "🥣"; 141 = "🥣"; 142 = "🫧"; 143 = "🫧"; 144 = "🥣"; 145 = "🥣";
    146 = "🌽"; 147 = "🍝"; 148 = "🌾"; 149 = "🌾"; 150 = "🍚"; 151 = "🍜";
    152 = "🌾"; 153 = "🌾"; 154 = "🌾"; 155 = "🍞"; 156 = "🍝"; 157 = "🍞";
    158 = "🍞"; 159 = "🍝"; 160 = "🌾"; 161 = "🌾"; 162 = "🥥"; 163 = "🫘";
    164 = "🌾"; 165 = "🌾"; 166 = "🌾"; 167 = "🌾"; 168 = "🥣"
}

foreach ($row in 134..168) {
    $ws.Range("F$row").Value = $emoji[$row]
}

$ws.Range("F168").Select()
